$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" column updates are plain numeric-looking strings (e.g. "306.68").
# Assigning those to .Value directly would make Excel auto-convert them to a
# Number, but the source data keeps them as plain text. Force text entry by
# switching those cells to the Text number format before writing the value,
# then restore the default "Normal" style so no stray formatting is left
# behind (matches cells such as D2/D14 whose new text still isn't numeric-
# looking and therefore never needed this treatment).
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '44.279.98'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '2.241.47'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').Value = '306.68'
$ws.Range('E5').Value = '  -2.86%  '
$ws.Range('D6').Value = '94.73'
$ws.Range('E6').Value = '  -4.69%  '
$ws.Range('D7').Value = '0.571'
$ws.Range('E7').Value = '  -0.76%  '
$ws.Range('E8').Value = '  +0.30%  '
$ws.Range('E9').Value = '  -1.64%  '
$ws.Range('D10').Value = '34.61'
$ws.Range('E10').Value = '  -4.77%  '
$ws.Range('D11').Value = '0.0812'
$ws.Range('E11').Value = '  -1.36%  '
$ws.Range('D12').Value = '7.18'
$ws.Range('E12').Value = '  -2.65%  '
$ws.Range('E13').Value = '  -0.31%  '
$ws.Range('D14').Value = '2.336.00'
$ws.Range('E14').Value = '  +3.97%  '
$ws.Range('D15').Value = '2.583.61'
$ws.Range('E15').Value = '  +0.00%  '
$ws.Range('D16').Value = '0.830'
$ws.Range('E16').Value = '  -1.77%  '
$ws.Range('E17').Value = '  -3.51%  '
$ws.Range('D18').Value = '44.045.24'
$ws.Range('E18').Value = '  +0.30%  '
$ws.Range('D19').Value = '0.0₃0966'
$ws.Range('E19').Value = '  -1.41%  '
$ws.Range('D20').Value = '6.40'
$ws.Range('E20').Value = '  +0.92%  '
$ws.Range('D21').Value = '12.19'
$ws.Range('E21').Value = '  -7.50%  '
$ws.Range('D22').Value = '65.47'
$ws.Range('E22').Value = '  -0.38%  '
$ws.Range('D23').Value = '237.60'
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('D24').Value = '2.95'
$ws.Range('E24').Value = '  -1.28%  '
$ws.Range('E25').Value = '  -1.13%  '
$ws.Range('D27').Value = '39.27'
$ws.Range('E27').Value = '  +7.57%  '
$ws.Range('D28').Value = '9.91'
$ws.Range('E28').Value = '  -2.45%  '
$ws.Range('E29').Value = '  +4.19%  '
$ws.Range('D30').Value = '20.03'
$ws.Range('E30').Value = '  -0.26%  '
$ws.Range('D31').Value = '5.84'
$ws.Range('E31').Value = '  -2.50%  '
$ws.Range('D32').Value = '153.30'
$ws.Range('E32').Value = '  -1.95%  '
$ws.Range('D33').Value = '0.0795'
$ws.Range('E33').Value = '  -5.48%  '
$ws.Range('E34').Value = '  -2.05%  '
$ws.Range('D35').Value = '3.14'
$ws.Range('E35').Value = '  -4.78%  '
$ws.Range('E36').Value = '  +1.68%  '
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('E38').Value = '  -7.44%  '
$ws.Range('E39').Value = '  -1.37%  '
$ws.Range('E40').Value = '  -4.75%  '
$ws.Range('D41').Value = '14.30'
$ws.Range('E41').Value = '  -7.30%  '
$ws.Range('E42').Value = '  -3.13%  '
$ws.Range('E43').Value = '  +0.34%  '
$ws.Range('D44').Value = '1.746.73'
$ws.Range('E44').Value = '  +2.50%  '
$ws.Range('D45').Value = '83.01'
$ws.Range('E45').Value = '  +0.31%  '
$ws.Range('D46').Value = '0.192'
$ws.Range('E46').Value = '  -1.54%  '
$ws.Range('D47').Value = '99.68'
$ws.Range('E47').Value = '  -2.18%  '
$ws.Range('E48').Value = '  -4.92%  '
$ws.Range('D49').Value = '1.61'
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').Value = '54.89'
$ws.Range('E50').Value = '  -2.80%  '
$ws.Range('B51').Value = 'FraxShare'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D51').Value = '8.08'
$ws.Range('E51').Value = '  -0.52%  '

$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D7').Style = "Normal"
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').Style = "Normal"
$ws.Range('D12').Style = "Normal"
$ws.Range('D16').Style = "Normal"
$ws.Range('D20').Style = "Normal"
$ws.Range('D21').Style = "Normal"
$ws.Range('D22').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D28').Style = "Normal"
$ws.Range('D30').Style = "Normal"
$ws.Range('D31').Style = "Normal"
$ws.Range('D32').Style = "Normal"
$ws.Range('D33').Style = "Normal"
$ws.Range('D35').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').Style = "Normal"
$ws.Range('D47').Style = "Normal"
$ws.Range('D49').Style = "Normal"
$ws.Range('D50').Style = "Normal"
$ws.Range('D51').Style = "Normal"
